$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing "Actual" value for week 2 (row 3). This drives the
# dependent Residual (H3) and Prediction-error (I3) formulas to recalc,
# clearing the previous #DIV/0! error.
$ws.Range("G3").Value = 6.61

# Add the next week's row (row 4): index, Week From, Week To.
$ws.Range("A4").Value = 3

# Copy the date formatting from the row above before writing the new dates
# so the new cells pick up the same date number format/style.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Value = 44885

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Value = 44892

# Clear the marching-ants clipboard marquee left by Copy().
$excel.CutCopyMode = 0

# Move the selection to reflect where the user left off entering data.
$ws.Range("D4").Select() | Out-Null
